# Updates the cryptos price/volume table to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain decimal number must keep their
# existing 'text' storage (e.g. trailing zero '575.40', not '575.4'), so force
# text formatting before writing the value, exactly like typing into a
# text-formatted cell in the Excel UI.
$textCells = @(
    'D5', 'D6', 'D9', 'D10', 'D12', 'D14', 'D18', 'D19', 'D21', 'D22', 'D23', 'D24', 'D25',
    'D26', 'D27', 'D29', 'D32', 'D33', 'D36', 'D37', 'D38', 'D40', 'D41', 'D42', 'D44', 'D45',
    'D46', 'D47', 'D49'
)
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the new cell values (Coin, Link, Price, Volume(1h) columns).
$ws.Range('D2').Value = '66.277.04'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '3.071.91'
$ws.Range('E3').Value = '  -1.60%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '575.40'
$ws.Range('E5').Value = '  -0.57%  '
$ws.Range('D6').Value = '170.62'
$ws.Range('E6').Value = '  -1.06%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '3.068.94'
$ws.Range('E8').Value = '  -1.55%  '
$ws.Range('D9').Value = '0.509'
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('D10').Value = '6.27'
$ws.Range('E10').Value = '  -2.01%  '
$ws.Range('E11').Value = '  -2.62%  '
$ws.Range('D12').Value = '0.469'
$ws.Range('E12').Value = '  -2.65%  '
$ws.Range('E13').Value = '  -3.85%  '
$ws.Range('D14').Value = '35.70'
$ws.Range('E14').Value = '  -3.88%  '
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').Value = '3.586.36'
$ws.Range('E16').Value = '  -1.41%  '
$ws.Range('D17').Value = '66.260.01'
$ws.Range('E17').Value = '  -1.05%  '
$ws.Range('D18').Value = '6.95'
$ws.Range('E18').Value = '  -2.91%  '
$ws.Range('D19').Value = '16.61'
$ws.Range('E19').Value = '  +2.28%  '
$ws.Range('D20').Value = '3.074.64'
$ws.Range('E20').Value = '  -1.45%  '
$ws.Range('D21').Value = '484.87'
$ws.Range('E21').Value = '  +2.02%  '
$ws.Range('D22').Value = '0.685'
$ws.Range('E22').Value = '  -3.39%  '
$ws.Range('D23').Value = '7.64'
$ws.Range('E23').Value = '  -2.91%  '
$ws.Range('D24').Value = '82.29'
$ws.Range('E24').Value = '  -1.85%  '
$ws.Range('D25').Value = '12.60'
$ws.Range('E25').Value = '  -4.83%  '
$ws.Range('D26').Value = '2.20'
$ws.Range('E26').Value = '  -3.39%  '
$ws.Range('D27').Value = '10.03'
$ws.Range('E27').Value = '  -2.86%  '
$ws.Range('E28').Value = '  -0.16%  '
$ws.Range('D29').Value = '7.83'
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('E30').Value = '  -5.21%  '
$ws.Range('E31').Value = '  -3.36%  '
$ws.Range('D32').Value = '27.63'
$ws.Range('E32').Value = '  -3.38%  '
$ws.Range('D33').Value = '0.111'
$ws.Range('E33').Value = '  -3.02%  '
$ws.Range('D34').Value = '0.0₃0914'
$ws.Range('E34').Value = '  -3.80%  '
$ws.Range('E35').Value = '  +0.06%  '
$ws.Range('D36').Value = '47.74'
$ws.Range('E36').Value = '  +1.77%  '
$ws.Range('D37').Value = '5.55'
$ws.Range('E37').Value = '  -5.01%  '
$ws.Range('D38').Value = '0.939'
$ws.Range('E38').Value = '  -4.01%  '
$ws.Range('E39').Value = '  -1.14%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '1.96'
$ws.Range('E40').Value = '  -4.76%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Value = '0.301'
$ws.Range('E41').Value = '  -3.50%  '
$ws.Range('D42').Value = '8.21'
$ws.Range('E42').Value = '  -4.57%  '
$ws.Range('D43').Value = '2.773.88'
$ws.Range('E43').Value = '  -1.74%  '
$ws.Range('D44').Value = '2.53'
$ws.Range('E44').Value = '  -0.73%  '
$ws.Range('D45').Value = '0.0343'
$ws.Range('E45').Value = '  -2.73%  '
$ws.Range('D46').Value = '134.66'
$ws.Range('E46').Value = '  -0.94%  '
$ws.Range('D47').Value = '363.83'
$ws.Range('E47').Value = '  -4.97%  '
$ws.Range('D49').Value = '24.21'
$ws.Range('E49').Value = '  -2.97%  '
$ws.Range('E50').Value = '  -2.59%  '
$ws.Range('E51').Value = '  -2.24%  '
